$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on cells whose new values look like numbers,
# so Excel keeps them as text (matching the original inlineStr cells)
# instead of silently converting them to floating point numbers.
$numericLookingRefs = @("D5", "D8", "D11", "D13", "D18", "D20", "D21", "D24", "D26", "D31", "D32", "D40", "D41", "D44", "D49")
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated price / volume values scraped by the GitHub Action run.
$ws.Range('D2').Value = '34.642.00'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '1.803.74'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '227.64'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').Value = '32.82'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '2.062.22'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').Value = '11.16'
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').Value = '1.795.63'
$ws.Range('E14').Value = '  +0.96%  '
$ws.Range('E15').Value = '  +2.75%  '
$ws.Range('D16').Value = '34.609.99'
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').Value = '68.94'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').Value = '0.0₃0806'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '247.50'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = '11.35'
$ws.Range('E21').Value = '  +3.42%  '
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('D24').Value = '170.14'
$ws.Range('E24').Value = '  +4.82%  '
$ws.Range('E25').Value = '  +2.06%  '
$ws.Range('D26').Value = '7.33'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  +11.06%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.0527'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.24'
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('D35').Value = '1.432.52'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('E36').Value = '  +8.35%  '
$ws.Range('E37').Value = '  +3.10%  '
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = '85.37'
$ws.Range('E40').Value = '  +6.48%  '
$ws.Range('D41').Value = '0.950'
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').Value = '13.85'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('E45').Value = '  +2.93%  '
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').Value = '1.961.35'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '105.89'
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('E51').Value = '  -4.62%  '

# Restore default number formatting / style on the cells we text-forced above
# so the workbook format table matches a normal (non text-forced) edit.
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "General"
    $ws.Range($ref).Style = "Normal"
}

Write-Output "Applied cryptos list update"